$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the indicator text (5.b.1) - new contact/metadata upload
$ws.Range("B4").Value = "5.b.1 Доля людей, имеющих мобильный телефон, в разбивке по полу "

# Update contact person
$ws.Range("B7").Value = "Калымбетова Ы.И."
$ws.Range("B7").Font.Name = "Calibri"

# Update contact phone
$ws.Range("B9").Value = "(0312) 32 46 55"
$ws.Range("B9").Font.Name = "Calibri"

# Update organization website
$ws.Range("B10").Value = "www.stat.gov.kg"
$ws.Range("B10").Font.Name = "Calibri"

# Update contact email
$ws.Range("B8").Value = "yryskan.kalymbetova@gmail.com"
$ws.Range("B8").Font.Name = "Calibri"

# Update organization / department
$ws.Range("B6").Value = "Национальный статистический комитет Кыргызской Республики" + [char]10 + "Управление статистики домашних хозяйств"
$ws.Range("B6").Font.Name = "Calibri"

# Reflect the saved selection state (active cell was B8 when the file was last saved)
$ws.Range("B8").Select()
